$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 30 is no longer the last row, so it should take on the regular
# data-row number format (matching A29 and the rows above it).
$ws.Range("A30").NumberFormat = $ws.Range("A29").NumberFormat

# Append the new daily entry (day 45771) as row 31, using the "last row"
# number format that row 30 used to have.
$ws.Range("A31").Value = 45771
$ws.Range("A31").NumberFormat = "YYYY-MM-DD"

$ws.Range("B31").Value = 124
$ws.Range("C31").Value = 126
$ws.Range("D31").Value = 125
